$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 8 new rows at the top; this shifts all previous rows (and the old
# shared H1:H5 formula) down by 8, landing the former row 1 at row 9.
$ws.Rows("1:8").Insert()

# The freshly inserted rows come back blank/unformatted - clone the number
# formats (date / text / text-left) from the row directly below (the old
# row 1, now row 9) onto the 8 new rows.
$ws.Range("A9:H9").Copy()
$ws.Range("A1:H8").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

# New transaction rows (most-recent-first), captured 2014-03-01.
$data = @(
    @(41703, "DB AH PROGRAMADO", "D", "0000948980", "AGENCIA PARA PROCESOS BATCH", "20.00  ", 567.02),
    @(41698, "INTERES A SU FAVOR", "C", "0003342171", "AGENCIA PARA PROCESOS BATCH", "0.02  ", 587.02),
    @(41698, "  TRANSFERENCIA INTERNET", "D", "0010403496", "AG. NORTE", "270.00  ", 587.00),
    @(41698, "13229623-MOVISTAR -CB-7872717", "D", "0008582680", "CENTRO DE SERVIC. OPERAT. SS.", "11.20  ", 857.00),
    @(41697, "  TRANSFERENCIA INTERNET", "C", "0000777219", "AG. NORTE", "115.00  ", 868.20),
    @(41696, "SPI COSTO OPER. CASH", "D", "0007678635", "CENTRO DE SERVIC. OPERAT. SS.", "0.27  ", 753.20),
    @(41696, "IVA COBRADO", "D", "0007678635", "CENTRO DE SERVIC. OPERAT. SS.", "0.03  ", 753.47),
    @(41696, "50-SPI-CCU - MIN. ECONOMIA CUENTA", "C", "0007678603", "TENA", "724.27  ", 753.50)
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $i + 1
    $row = $data[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
    $ws.Cells.Item($r, 7).Value = $row[6]

    # Same CONCATENATE(...) builder formula as every other row, anchored to
    # this row so relative refs (A<r>, B<r>, ...) resolve correctly.
    $ws.Cells.Item($r, 8).Formula = "=CONCATENATE(""array('mo_fecha' => new \DateTime('"",TEXT(A$r,""yyyy-mm-dd""),""'), 'mo_concepto' => '"",B$r,""', 'mo_tipo' => '"",C$r,""', 'mo_documento' => '"",D$r,""', 'mo_oficina' => '"",E$r,""', 'mo_monto' => "",F$r,"", 'mo_saldo' => "",G$r,"", 'mo_fecha_crea' => new \DateTime('"",TEXT(NOW(),""yyyy-mm-dd H:m:s""),""'), 'mo_quien_crea' => 1, 'mo_fecha_modifica' => NULL, 'mo_quien_modifica' => NULL, 'mo_borrado_logico' => false),"")"
}

$ws.Range("H1").Select()
$ws.Range("H1:H8").Select()
